$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.285.28'
$ws.Range("E2").Value = '  -3.98%  '

$ws.Range("D3").Value = '2.980.35'
$ws.Range("E3").Value = '  -3.65%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.06%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '2.974.83'
$ws.Range("E8").Value = '  -3.60%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.494'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.147'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.08%  '

$ws.Range("E12").Value = '  -3.74%  '

$ws.Range("E13").Value = '  -3.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.03%  '

$ws.Range("D15").Value = '3.468.68'
$ws.Range("E15").Value = '  -3.47%  '

$ws.Range("E16").Value = '  -2.00%  '

$ws.Range("D17").Value = '61.356.92'
$ws.Range("E17").Value = '  -3.88%  '

$ws.Range("D18").Value = '2.985.96'
$ws.Range("E18").Value = '  -3.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.98%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.668'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.48%  '

$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("E30").Value = '  +2.58%  '

$ws.Range("E31").Value = '  -2.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.42'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '54.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.86%  '

$ws.Range("B34").Value = 'Stacks'
$ws.Range("C34").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.26'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.32%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.23%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.85'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.82%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '451.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -10.18%  '

$ws.Range("D38").Value = '3.145.18'
$ws.Range("E38").Value = '  -4.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0782'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0381'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.11%  '

$ws.Range("E41").Value = '  +0.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.70%  '

$ws.Range("E43").Value = '  -9.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.43%  '

$ws.Range("E45").Value = '  +0.10%  '

$ws.Range("E46").Value = '  -7.36%  '

$ws.Range("E47").Value = '  -6.25%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.107'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.23%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '117.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.83%  '

$ws.Range("E50").Value = '  -9.46%  '

$ws.Range("E51").Value = '  +6.39%  '
